# Update "想去人数" (interested-count) figures in F column across sheets,
# matching regenerated output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(3, 6).Value  = 1004
$ws1.Cells.Item(6, 6).Value  = 715
$ws1.Cells.Item(7, 6).Value  = 251
$ws1.Cells.Item(9, 6).Value  = 32
$ws1.Cells.Item(12, 6).Value = 76
$ws1.Cells.Item(13, 6).Value = 837
$ws1.Cells.Item(15, 6).Value = 1984
$ws1.Cells.Item(16, 6).Value = 479
$ws1.Cells.Item(17, 6).Value = 7145
$ws1.Cells.Item(20, 6).Value = 54
$ws1.Cells.Item(21, 6).Value = 91
$ws1.Cells.Item(22, 6).Value = 17
$ws1.Cells.Item(23, 6).Value = 219

# Sheet "本地生活" (local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2, 6).Value = 5474
$ws3.Cells.Item(3, 6).Value = 390
$ws3.Cells.Item(4, 6).Value = 382

# Sheet "全部类型" (all types combined)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(3, 6).Value  = 5474
$ws4.Cells.Item(4, 6).Value  = 390
$ws4.Cells.Item(5, 6).Value  = 382
$ws4.Cells.Item(7, 6).Value  = 1004
$ws4.Cells.Item(12, 6).Value = 715
$ws4.Cells.Item(13, 6).Value = 251
$ws4.Cells.Item(16, 6).Value = 32
$ws4.Cells.Item(20, 6).Value = 76
$ws4.Cells.Item(22, 6).Value = 837
$ws4.Cells.Item(25, 6).Value = 1984
$ws4.Cells.Item(26, 6).Value = 479
$ws4.Cells.Item(27, 6).Value = 7145
$ws4.Cells.Item(31, 6).Value = 54
$ws4.Cells.Item(32, 6).Value = 91
$ws4.Cells.Item(34, 6).Value = 17
$ws4.Cells.Item(35, 6).Value = 219
